$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update column F (想去人数) for rows 2-9
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 708
$ws1.Range("F3").Value = 34
$ws1.Range("F4").Value = 239
$ws1.Range("F5").Value = 2403
$ws1.Range("F6").Value = 50
$ws1.Range("F7").Value = 3546
$ws1.Range("F8").Value = 465
$ws1.Range("F9").Value = 903

# Sheet "全部类型" (fourth sheet) - update column F (想去人数) for rows 2-10
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 708
$ws4.Range("F3").Value = 34
$ws4.Range("F5").Value = 239
$ws4.Range("F6").Value = 2404
$ws4.Range("F7").Value = 50
$ws4.Range("F8").Value = 3546
$ws4.Range("F9").Value = 465
$ws4.Range("F10").Value = 903
